$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 386, shifting existing rows 386-413 down to 387-414.
$ws.Rows.Item(386).Insert()

# Populate the newly inserted row 386 with the new record's data.
$ws.Cells.Item(386, 1).Value = 6
$ws.Cells.Item(386, 2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(386, 3).Value = 'Metropolitana'
$ws.Cells.Item(386, 4).Value = 44746
$ws.Cells.Item(386, 5).Value = 13
$ws.Cells.Item(386, 6).Value = 100112032
$ws.Cells.Item(386, 7).Value = 'Zapallo italiano'
$ws.Cells.Item(386, 8).Value = 'Sin especificar'
$ws.Cells.Item(386, 9).Value = 'Primera'
$ws.Cells.Item(386, 10).Value = 400
$ws.Cells.Item(386, 11).Value = 10000
$ws.Cells.Item(386, 12).Value = 12000
$ws.Cells.Item(386, 13).Value = 10850
$ws.Cells.Item(386, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(386, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(386, 16).Value = 217
$ws.Cells.Item(386, 17).Value = 50
$ws.Cells.Item(386, 18).Value = 'Hortaliza'

# Ensure the date cell keeps the existing date-number-format style (s="2"),
# matching the style used by the other date cells in column D.
$ws.Cells.Item(386, 4).NumberFormat = $ws.Cells.Item(387, 4).NumberFormat
